$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = 'Normal'
}

Set-TextValue 'D2' '25.912.64'
Set-TextValue 'E2' '  -1.30%  '
Set-TextValue 'D3' '1.638.19'
Set-TextValue 'E3' '  -0.42%  '
Set-TextValue 'E4' '  +0.39%  '
Set-TextValue 'D5' '214.78'
Set-TextValue 'E5' '  -0.91%  '
Set-TextValue 'D6' '0.506'
Set-TextValue 'E6' '  -0.08%  '
Set-TextValue 'E7' '  +0.39%  '
Set-TextValue 'E8' '  -1.38%  '
Set-TextValue 'D9' '0.0638'
Set-TextValue 'E9' '  +0.00%  '
Set-TextValue 'D10' '19.60'
Set-TextValue 'E10' '  -2.15%  '
Set-TextValue 'D11' '0.0794'
Set-TextValue 'E11' '  +0.09%  '
Set-TextValue 'D12' '1.865.75'
Set-TextValue 'E12' '  -0.31%  '
Set-TextValue 'B13' 'WrappedEther'
Set-TextValue 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.665.87'
Set-TextValue 'E13' '  +2.52%  '
Set-TextValue 'B14' 'Polkadot'
Set-TextValue 'C14' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D14' '4.25'
Set-TextValue 'E14' '  -1.10%  '
Set-TextValue 'D15' '0.543'
Set-TextValue 'E15' '  -2.20%  '
Set-TextValue 'E16' '  -0.83%  '
Set-TextValue 'D17' '62.67'
Set-TextValue 'E17' '  -1.54%  '
Set-TextValue 'D18' '25.928.65'
Set-TextValue 'E18' '  -1.08%  '
Set-TextValue 'E19' '  +0.33%  '
Set-TextValue 'D20' '193.46'
Set-TextValue 'E20' '  +0.02%  '
Set-TextValue 'E21' '  -1.94%  '
Set-TextValue 'D22' '9.91'
Set-TextValue 'E22' '  -1.80%  '
Set-TextValue 'D23' '6.27'
Set-TextValue 'E23' '  -1.28%  '
Set-TextValue 'D24' '143.83'
Set-TextValue 'E24' '  +0.34%  '
Set-TextValue 'B25' 'BinanceUSD'
Set-TextValue 'C25' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D25' '1.00'
Set-TextValue 'E25' '  +0.51%  '
Set-TextValue 'B26' 'Toncoin'
Set-TextValue 'C26' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D26' '1.78'
Set-TextValue 'E26' '  -0.69%  '
Set-TextValue 'D27' '0.127'
Set-TextValue 'E27' '  +1.29%  '
Set-TextValue 'D28' '6.83'
Set-TextValue 'E28' '  -1.77%  '
Set-TextValue 'D29' '15.49'
Set-TextValue 'E29' '  -0.91%  '
Set-TextValue 'E30' '  -0.77%  '
Set-TextValue 'D31' '0.0503'
Set-TextValue 'E31' '  +0.13%  '
Set-TextValue 'E32' '  -2.00%  '
Set-TextValue 'E33' '  -0.96%  '
Set-TextValue 'D34' '1.53'
Set-TextValue 'E34' '  -4.27%  '
Set-TextValue 'E35' '  +1.23%  '
Set-TextValue 'D36' '0.902'
Set-TextValue 'E36' '  -1.65%  '
Set-TextValue 'D37' '1.136.44'
Set-TextValue 'E37' '  -0.18%  '
Set-TextValue 'D38' '0.543'
Set-TextValue 'E38' '  -2.13%  '
Set-TextValue 'D39' '2.46'
Set-TextValue 'E39' '  -2.03%  '
Set-TextValue 'E40' '  -0.63%  '
Set-TextValue 'B41' 'PaxDollar'
Set-TextValue 'C41' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D41' '1.00'
Set-TextValue 'E41' '  +0.33%  '
Set-TextValue 'D42' '0.806'
Set-TextValue 'E42' '  +1.00%  '
Set-TextValue 'B43' 'Quant'
Set-TextValue 'C43' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D43' '99.48'
Set-TextValue 'E43' '  -1.03%  '
Set-TextValue 'B44' 'FraxShare'
Set-TextValue 'C44' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D44' '5.44'
Set-TextValue 'E44' '  -3.36%  '
Set-TextValue 'B45' 'RocketPoolETH'
Set-TextValue 'C45' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 'D45' '1.775.30'
Set-TextValue 'E45' '  -0.29%  '
Set-TextValue 'B46' 'BabyDogeCoin'
Set-TextValue 'C46' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D46' '0.0₆0116'
Set-TextValue 'E46' '  +14.86%  '
Set-TextValue 'B47' 'Aave'
Set-TextValue 'C47' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D47' '56.57'
Set-TextValue 'E47' '  +0.63%  '
Set-TextValue 'B48' 'Cronos'
Set-TextValue 'C48' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D48' '0.0532'
Set-TextValue 'E48' '  +3.01%  '
Set-TextValue 'B49' 'RenderToken'
Set-TextValue 'C49' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D49' '1.46'
Set-TextValue 'E49' '  -2.39%  '
Set-TextValue 'B50' 'EnergySwap'
Set-TextValue 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '7.65'
Set-TextValue 'E50' '  -1.27%  '
Set-TextValue 'B51' 'Mantle'
Set-TextValue 'C51' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D51' '0.415'
Set-TextValue 'E51' '  -0.75%  '
